$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'report_id'
$ws.Range("C2").Value = 'number'
$ws.Range("D2").Value = '7-digit ICSR number (no spaces).'

$ws.Range("B3").Value = 'submission_date'
$ws.Range("C3").Value = 'string'
$ws.Range("D3").Value = 'Date report was received by CTP; this is the earliest date of report receipt, either to Safety Reporting Portal (SRP) or by other means'
$ws.Rows.Item(3).RowHeight = 34

$ws.Range("B4").Value = 'number_tobacco_products'
$ws.Range("C4").Value = 'number'
$ws.Range("D4").Value = 'System-calculated number of Tobacco Product Problems reported, displayed as a whole number, ≥ 0.'
$ws.Rows.Item(4).RowHeight = 34

$ws.Range("B5").Value = 'tobacco_products'
$ws.Range("C5").Value = 'array'
$ws.Range("D5").Value = 'Text reflecting the SRP tobacco Product Type selected by the reporter.'
$ws.Rows.Item(5).RowHeight = 17

$ws.Range("B6").Value = 'number_health_problems'
$ws.Range("C6").Value = 'number'
$ws.Range("D6").Value = 'System-calculated number of Health Problems (i.e., MedDRA terms selected from a standardized list of symptoms, signs, diagnoses and outcomes) reported, displayed as a whole number, ≥0.'
$ws.Rows.Item(6).RowHeight = 51

$ws.Range("B7").Value = 'reported_health_problems'
$ws.Range("C7").Value = 'array'
$ws.Range("D7").Value = 'Text reflecting the MedDRA terms selected by the reporter.'
$ws.Rows.Item(7).RowHeight = 17

$ws.Range("B8").Value = 'nonuser_affected'
$ws.Range("C8").Value = 'string'
$ws.Range("D8").Value = 'Displays text reflecting the response to this optional question (2017 - 12/14/2018) or required question (12/15/2018 onward) as “No information provided” if not answered, or Yes/No.'
$ws.Rows.Item(8).RowHeight = 51

$ws.Range("B9").Value = 'number_product_problems'
$ws.Range("C9").Value = 'number'
$ws.Range("D9").Value = 'System-calculated number of categorical Product Problems reported, displayed as a whole number, ≥0.'
$ws.Rows.Item(9).RowHeight = 34

$ws.Range("B10").Value = 'reported_product_problems'
$ws.Range("C10").Value = 'array'
$ws.Range("D10").Value = 'Text reflecting the SRP categorical list of values.'
$ws.Rows.Item(10).RowHeight = 17

$ws.Range("B12").Select()

$win = $wb.Windows.Item(1)
$win.Left = 11660
$win.Top = 460